$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.6753301551942219, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 2.997429241610044)
    3 = @(0.003994804209775715, 0.002777888934908601, 26.21740644021617, 8.660232485948974, 34.88441161930983)
    4 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    5 = @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 5.553084769722144)
    6 = @(3.230985683306322, 0.3127903958511391, 0.8054896365839992, 0.496779210170732, 4.846044925912192)
    7 = @(3.230985683306322, 3099.503889238888, 0.8054896365839992, 8.660232485948974, 3112.200597044728)
    8 = @(0.6753301551942219, 1.667794583268128, 0.1575252929769615, 0.496779210170732, 2.997429241610044)
    9 = @(0.127881588408715, 0.3127903958511391, 0.8054896365839992, 0.496779210170732, 1.742940831014585)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
